# Failsafe pulse generator adapted for different requirements of WS2813B-V5:
# append the new 555-oscillator measurement rows and drop the now-unused
# helper column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("4" in G2) is no longer part of the table -> clear it so the
# sheet's used range shrinks back down to columns A:F.
$ws.Range("G2").ClearContents()

# New measurement rows (C, RA, RB, t_h, t_l, T) appended below the existing
# 8 rows of data, for a WS2813B-V5 oscillator.
$capVal = 0.00000000047
$data = @(
    @($capVal, 1500, 1000, 1315, 641.7,  1957),
    @($capVal, 1500, 1500, 1479, 857.4,  2336),
    @($capVal, 1500, 1300, 1406, 768,    2174),
    @($capVal, 1000, 1300, 1190, 787.8,  1977),
    @($capVal,  910, 1300, 1155, 794.4,  1949),
    @($capVal,  680, 1300, 1050, 821,    1871),
    @(470,      680, 1200, 1004, 767,    1771),
    @($capVal,  750, 1200, 1053, 776.3,  1830),
    @($capVal,  820, 1200, 1070, 762.1,  1832),
    @($capVal,  910, 1200, 1130, 771.2,  1901)
)

$row = 9
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

# The capacitance column (A) uses a scientific-notation number format;
# carry it down onto the freshly-written rows, matching A2:A8.
$ws.Range("A9:A18").NumberFormat = $ws.Range("A2").NumberFormat

# Leave the selection where Excel would land after typing the last row.
[void]$ws.Range("A19").Select()
